# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 1
    3  = 1
    4  = 1
    5  = 1
    6  = 2
    7  = 1
    8  = 0
    9  = 3
    10 = 0
    11 = 0
    12 = 6
    13 = 3
    14 = 1
    15 = 3
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
